$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @(2026, 2, 24, 4659, 4476, 4406, 4405, 4447, 4603, 5128, 5700, 6154, 6281, 6142, 6035, 5921, 5835, 5720, 5700, 5715, 5919, 6400, 6425, 6234, 5992, 5770, 5429)

$row = 21
for ($i = 0; $i -lt $values.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item($row, $col).Value = $values[$i]
}
